# New column: alamat, pdf file fixes
#
# Insert a new "alamat" column (header: "alamatku") right before the
# existing "namaayah" column (old column E), shifting namaayah, namaibu,
# haribaptis, the baptism-date column and the pastor column one place to
# the right. Then remove the trailing PDF filename column (old column J,
# e.g. "001_GPT_A_V_19") which is no longer needed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new blank column at E (pushes E:J -> F:K)
$ws.Range("E1").EntireColumn.Insert()

# 2) Make sure the new column has no inherited number formatting/style
#    (it would otherwise pick up the date style from column D)
$ws.Range("E1:E10").ClearFormats()

# 3) Fill the new "alamat" column with its value for every data row
$ws.Range("E1:E10").Value = "alamatku"

# 4) Remove the old PDF-filename column, which is now column K
$ws.Range("K1").EntireColumn.Delete()

# 5) Reflect the final selection used when the workbook was last saved
$ws.Range("K1:K10").Select()
